# Update the "precision" (column B) and "f1-score" (column C) values in the
# DecisionTree classification-report worksheet to reflect the refreshed
# model evaluation numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3; Precision = 0.4252026529108327; F1Score = 0.4333458505444986 }
    @{ Row = 5; Precision = 0.5649289099526066; F1Score = 0.5635933806146571 }
    @{ Row = 6; Precision = 0.1016597510373444; F1Score = 0.09099350046425254 }
    @{ Row = 7; Precision = 0.04330708661417323; F1Score = 0.0352 }
    @{ Row = 8; Precision = 0.02649006622516556; F1Score = 0.02228412256267409 }
    @{ Row = 9; Precision = 0; F1Score = 0 }
    @{ Row = 10; Precision = 0.06179775280898876; F1Score = 0.04988662131519275 }
    @{ Row = 11; Precision = 0.09508196721311475; F1Score = 0.0703883495145631 }
    @{ Row = 12; Precision = 0.157819225251076; F1Score = 0.1403956604977664 }
    @{ Row = 13; Precision = 0.1062992125984252; F1Score = 0.08925619834710743 }
    @{ Row = 14; Precision = 0.08; F1Score = 0.05825242718446601 }
    @{ Row = 15; Precision = 0.02362204724409449; F1Score = 0.01973684210526316 }
    @{ Row = 16; Precision = 0; F1Score = 0 }
    @{ Row = 17; Precision = 0.05161290322580645; F1Score = 0.03950617283950617 }
    @{ Row = 18; Precision = 0.07906976744186046; F1Score = 0.06261510128913443 }
    @{ Row = 19; Precision = 0.1734390485629336; F1Score = 0.1769464105156724 }
    @{ Row = 20; Precision = 0.08823529411764706; F1Score = 0.072992700729927 }
    @{ Row = 21; Precision = 0.05990783410138249; F1Score = 0.04585537918871253 }
    @{ Row = 22; Precision = 0.04761904761904762; F1Score = 0.0472636815920398 }
    @{ Row = 23; Precision = 0.04081632653061224; F1Score = 0.032 }
    @{ Row = 24; Precision = 0; F1Score = 0 }
    @{ Row = 25; Precision = 0.01639344262295082; F1Score = 0.01574803149606299 }
    @{ Row = 28; Precision = 0.0622568093385214; F1Score = 0.05555555555555555 }
    @{ Row = 29; Precision = 0.4440742503569728; F1Score = 0.4356759280877889 }
    @{ Row = 30; Precision = 0.1491228070175439; F1Score = 0.1392301392301392 }
    @{ Row = 31; Precision = 0.2335907335907336; F1Score = 0.1917591125198098 }
    @{ Row = 32; Precision = 0; F1Score = 0 }
    @{ Row = 33; Precision = 0.3086876155268022; F1Score = 0.2661354581673307 }
    @{ Row = 34; Precision = 0.03260869565217391; F1Score = 0.02553191489361702 }
    @{ Row = 35; Precision = 0.07364341085271318; F1Score = 0.0581039755351682 }
    @{ Row = 36; Precision = 0.3824884792626728; F1Score = 0.3908509922637067 }
)

foreach ($update in $updates) {
    $ws.Cells.Item($update.Row, 2).Value = $update.Precision
    $ws.Cells.Item($update.Row, 3).Value = $update.F1Score
}

